$d = $word.ActiveDocument

$replacements = @(
    @{old="50×77="; new="89×25="},
    @{old="60×93="; new="13×12="},
    @{old="70×60="; new="34×16="},
    @{old="69×48="; new="89×73="},
    @{old="95×66="; new="68×22="},
    @{old="33×17="; new="21×92="},
    @{old="43×56="; new="68×31="},
    @{old="15×37="; new="22×20="},
    @{old="59×59="; new="39×26="},
    @{old="58×80="; new="15×38="},
    @{old="30×42="; new="80×49="},
    @{old="72×34="; new="20×86="},
    @{old="14×80="; new="80×48="},
    @{old="64×14="; new="66×31="},
    @{old="13×55="; new="43×72="},
    @{old="90×77="; new="81×46="},
    @{old="17×91="; new="87×14="},
    @{old="98×40="; new="47×91="},
    @{old="34×20="; new="84×90="},
    @{old="91×46="; new="72×76="},
    @{old="64×86="; new="46×56="},
    @{old="33×48="; new="21×55="},
    @{old="88×94="; new="26×12="},
    @{old="74×86="; new="25×54="},
    @{old="81×86="; new="24×73="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}

$d.Save()
